$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same (row 1): NumeroTransaccion, Referencia, Resultado, Mensaje, TiempoTransaccion

# Row 2
$ws.Cells.Item(2, 1).Value = "'1"
$ws.Cells.Item(2, 2).Value = "00837202409271520249-01102024"
$ws.Cells.Item(2, 3).Value = "BE"
$ws.Cells.Item(2, 4).Value = "Error BusinessRule process state: No se han encontrado casos pendientes de procesar"
$ws.Cells.Item(2, 5).Value = "TransactionTime: 0h 2m 2s"

# Row 3
$ws.Cells.Item(3, 1).Value = "'2"
$ws.Cells.Item(3, 2).Value = "00837202409271514903-01102024"
$ws.Cells.Item(3, 3).Value = "OK"
$ws.Cells.Item(3, 4).Value = "Successful Transaction: 2 - Reference: 00837202409271514903-01102024"
$ws.Cells.Item(3, 5).Value = "TransactionTime: 0h 0m 34s"

# Row 4
$ws.Cells.Item(4, 1).Value = "'3"
$ws.Cells.Item(4, 2).Value = "00837202409271514879-01102024"
$ws.Cells.Item(4, 3).Value = "BE"
$ws.Cells.Item(4, 4).Value = "Error BusinessRule process state: No se han encontrado casos pendientes de procesar"
$ws.Cells.Item(4, 5).Value = "TransactionTime: 0h 1m 23s"

# Row 5
$ws.Cells.Item(5, 1).Value = "'4"
$ws.Cells.Item(5, 2).Value = "00837202409271514149-01102024"
$ws.Cells.Item(5, 3).Value = "BE"
$ws.Cells.Item(5, 4).Value = "Error BusinessRule process state: No se han encontrado casos pendientes de procesar"
$ws.Cells.Item(5, 5).Value = "TransactionTime: 0h 1m 0s"

# Row 6
$ws.Cells.Item(6, 1).Value = "'5"
$ws.Cells.Item(6, 2).Value = "00837202409271513999-01102024"
$ws.Cells.Item(6, 3).Value = "BE"
$ws.Cells.Item(6, 4).Value = "Error BusinessRule process state: No se han encontrado casos pendientes de procesar"
$ws.Cells.Item(6, 5).Value = "TransactionTime: 0h 0m 59s"

# Row 7
$ws.Cells.Item(7, 1).Value = "'6"
$ws.Cells.Item(7, 2).Value = "00837202409161511213-18092024"
$ws.Cells.Item(7, 3).Value = "BE"
$ws.Cells.Item(7, 4).Value = "Error BusinessRule process state: No se encontrarón datos en detalle del item: 00837202409161511213"
$ws.Cells.Item(7, 5).Value = "TransactionTime: 0h 0m 42s"
